$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.436534333333333
$ws.Cells.Item(2, 8).Value = 4.309603
$ws.Cells.Item(2, 9).Value = 0.03241561610838976
$ws.Cells.Item(2, 10).Value = 0.03241561610838976
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.09554499999999999
$ws.Cells.Item(2, 14).Value = 0.286635
$ws.Cells.Item(2, 15).Value = 0.01029975823317688
$ws.Cells.Item(2, 16).Value = 0.01029975823317688
$ws.Cells.Item(2, 17).Value = 0.1372536728783333
$ws.Cells.Item(2, 18).Value = 1.235283055905
$ws.Cells.Item(2, 19).Value = 0.0003338730088958885
$ws.Cells.Item(2, 20).Value = 0.0003338730088958885

$ws.Cells.Item(3, 7).Value = 1.436534333333333
$ws.Cells.Item(3, 8).Value = 4.309603
$ws.Cells.Item(3, 9).Value = 0.03241561610838976
$ws.Cells.Item(3, 10).Value = 0.03241561610838976
$ws.Cells.Item(3, 15).Value = 0.4011437372432085
$ws.Cells.Item(3, 16).Value = 0.4011437372432086
$ws.Cells.Item(3, 17).Value = 5.345606182426778
$ws.Cells.Item(3, 18).Value = 48.110455641841
$ws.Cells.Item(3, 19).Value = 0.01300332139076062
$ws.Cells.Item(3, 20).Value = 0.01300332139076062

$ws.Cells.Item(4, 7).Value = 1.436534333333333
$ws.Cells.Item(4, 8).Value = 4.309603
$ws.Cells.Item(4, 9).Value = 0.03241561610838976
$ws.Cells.Item(4, 10).Value = 0.03241561610838976
$ws.Cells.Item(4, 13).Value = 5.459703999999999
$ws.Cells.Item(4, 14).Value = 16.379112
$ws.Cells.Item(4, 15).Value = 0.5885565045236145
$ws.Cells.Item(4, 16).Value = 0.5885565045236146
$ws.Cells.Item(4, 17).Value = 7.843052245837333
$ws.Cells.Item(4, 18).Value = 70.587470212536
$ws.Cells.Item(4, 19).Value = 0.01907842170873325
$ws.Cells.Item(4, 20).Value = 0.01907842170873325

$ws.Cells.Item(5, 7).Value = 34.88211266666666
$ws.Cells.Item(5, 9).Value = 0.7871201871162607
$ws.Cells.Item(5, 10).Value = 0.7871201871162609
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.09554499999999999
$ws.Cells.Item(5, 14).Value = 0.286635
$ws.Cells.Item(5, 15).Value = 0.01029975823317688
$ws.Cells.Item(5, 16).Value = 0.01029975823317688
$ws.Cells.Item(5, 17).Value = 3.332811454736666
$ws.Cells.Item(5, 18).Value = 29.99530309263
$ws.Cells.Item(5, 19).Value = 0.008107147627750434
$ws.Cells.Item(5, 20).Value = 0.008107147627750435

$ws.Cells.Item(6, 7).Value = 34.88211266666666
$ws.Cells.Item(6, 9).Value = 0.7871201871162607
$ws.Cells.Item(6, 10).Value = 0.7871201871162609
$ws.Cells.Item(6, 15).Value = 0.4011437372432085
$ws.Cells.Item(6, 16).Value = 0.4011437372432086
$ws.Cells.Item(6, 19).Value = 0.3157483335193904
$ws.Cells.Item(6, 20).Value = 0.3157483335193905

$ws.Cells.Item(7, 7).Value = 34.88211266666666
$ws.Cells.Item(7, 9).Value = 0.7871201871162607
$ws.Cells.Item(7, 10).Value = 0.7871201871162609
$ws.Cells.Item(7, 13).Value = 5.459703999999999
$ws.Cells.Item(7, 14).Value = 16.379112
$ws.Cells.Item(7, 15).Value = 0.5885565045236145
$ws.Cells.Item(7, 16).Value = 0.5885565045236146
$ws.Cells.Item(7, 17).Value = 190.4460100546506
$ws.Cells.Item(7, 18).Value = 1714.014090491856
$ws.Cells.Item(7, 19).Value = 0.4632647059691198
$ws.Cells.Item(7, 20).Value = 0.46326470596912

$ws.Cells.Item(8, 7).Value = 7.997472999999999
$ws.Cells.Item(8, 8).Value = 23.992419
$ws.Cells.Item(8, 9).Value = 0.1804641967753495
$ws.Cells.Item(8, 10).Value = 0.1804641967753495
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.09554499999999999
$ws.Cells.Item(8, 14).Value = 0.286635
$ws.Cells.Item(8, 15).Value = 0.01029975823317688
$ws.Cells.Item(8, 16).Value = 0.01029975823317688
$ws.Cells.Item(8, 17).Value = 0.7641185577849998
$ws.Cells.Item(8, 18).Value = 6.877067020064999
$ws.Cells.Item(8, 19).Value = 0.001858737596530559
$ws.Cells.Item(8, 20).Value = 0.001858737596530559

$ws.Cells.Item(9, 7).Value = 7.997472999999999
$ws.Cells.Item(9, 8).Value = 23.992419
$ws.Cells.Item(9, 9).Value = 0.1804641967753495
$ws.Cells.Item(9, 10).Value = 0.1804641967753495
$ws.Cells.Item(9, 15).Value = 0.4011437372432085
$ws.Cells.Item(9, 16).Value = 0.4011437372432086
$ws.Cells.Item(9, 17).Value = 29.76005523891033
$ws.Cells.Item(9, 18).Value = 267.840497150193
$ws.Cells.Item(9, 19).Value = 0.07239208233305747
$ws.Cells.Item(9, 20).Value = 0.07239208233305748

$ws.Cells.Item(10, 7).Value = 7.997472999999999
$ws.Cells.Item(10, 8).Value = 23.992419
$ws.Cells.Item(10, 9).Value = 0.1804641967753495
$ws.Cells.Item(10, 10).Value = 0.1804641967753495
$ws.Cells.Item(10, 13).Value = 5.459703999999999
$ws.Cells.Item(10, 14).Value = 16.379112
$ws.Cells.Item(10, 15).Value = 0.5885565045236145
$ws.Cells.Item(10, 16).Value = 0.5885565045236146
$ws.Cells.Item(10, 17).Value = 43.66383532799199
$ws.Cells.Item(10, 18).Value = 392.974517951928
$ws.Cells.Item(10, 19).Value = 0.1062133768457614
$ws.Cells.Item(10, 20).Value = 0.1062133768457615
